# Update the AC cable reactance matrix (symmetric cluster-to-cluster values)
# on Sheet1. Only the numeric body (B2:I9) changes; the header row/column
# labels stay the same.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Rows 2-9 correspond to cluster_0..cluster_7 (column A labels),
# columns B-I (2-9) correspond to the same clusters.
$values = @(
    @(0,   13, 6,  0,   1.2, 24, 0,  12),
    @(13,  0,  4,  0,   0,   23, 8,  0),
    @(6,   4,  0,  0,   0,   0,  0,  23),
    @(0,   0,  0,  0,   0,   4.5,37, 0),
    @(1.2, 0,  0,  0,   0,   20, 0,  0),
    @(24,  23, 0,  4.5, 20,  0,  0,  0),
    @(0,   8,  0,  37,  0,   0,  0,  0),
    @(12,  0,  23, 0,   0,   0,  0,  0)
)

for ($r = 0; $r -lt 8; $r++) {
    $row = $values[$r]
    for ($c = 0; $c -lt 8; $c++) {
        $ws.Cells.Item($r + 2, $c + 2).Value = $row[$c]
    }
}
